$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set columns C and D (date / time) to Text format BEFORE writing values so
# that date/time-looking strings are kept as literal text instead of being
# auto-converted to Excel date/time serial numbers.
$ws.Range("C1:D3").NumberFormat = "@"

# Header row
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "event_type"
$ws.Range("C1").Value = "date"
$ws.Range("D1").Value = "time"
$ws.Range("E1").Value = "sponsor_name"
$ws.Range("F1").Value = "participation"

# Row 2
$ws.Range("A2").Value = "desayuno con COLGATE"
$ws.Range("B2").Value = "F"
$ws.Range("C2").Value = "2023-10-03"
$ws.Range("D2").Value = "08:30:00"
$ws.Range("E2").Value = "COLGATE"
$ws.Range("F2").Value = "desayuno con representate de colgate"

# Row 3
$ws.Range("A3").Value = "Reunion de control"
$ws.Range("B3").Value = "O"
$ws.Range("C3").Value = "2024-01-07"
$ws.Range("D3").Value = "16:45:00"
$ws.Range("E3").Value = "LAURA MEDINA"
$ws.Range("F3").Value = "control de calidad"

# Column widths (best-fit sized to the content, matching the authoring app)
$ws.Columns.Item(1).ColumnWidth = 21.85546875
$ws.Columns.Item(2).ColumnWidth = 11.140625
$ws.Columns.Item(3).ColumnWidth = 10.7109375
$ws.Columns.Item(4).ColumnWidth = 8.140625
$ws.Columns.Item(5).ColumnWidth = 14.5703125
$ws.Columns.Item(6).ColumnWidth = 35.28515625

# Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Restore selection to the cell below/right of the inserted table
$ws.Range("D4").Select() | Out-Null
